$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New backlog item for the "Mejoras" (improvements) list: add the text
# "Incluir precio de adornos, canasta, etc" to column D ("Mejora Funcional")
# on the first free row (row 10).
$ws.Range("D10").Value = "Incluir precio de adornos, canasta, etc"

# Excel persists the last selected cell in the sheet view; after typing the
# new entry the cursor moved down to the next empty row.
$ws.Range("D11").Select()
